$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting existing rows 5..102 down to 6..103.
$ws.Rows(5).Insert()

# Populate the newly inserted row 5 with the new record's data.
$ws.Cells.Item(5, 1).Value = 7
$ws.Cells.Item(5, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value = "Ñuble"
$ws.Cells.Item(5, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 100112045
$ws.Cells.Item(5, 7).Value = "Zapallo"
$ws.Cells.Item(5, 8).Value = "Paine"
$ws.Cells.Item(5, 9).Value = "1a (guarda)"
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 220
$ws.Cells.Item(5, 12).Value = 250
$ws.Cells.Item(5, 13).Value = 235
$ws.Cells.Item(5, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(5, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(5, 16).Value = 235
$ws.Cells.Item(5, 17).Value = 1
$ws.Cells.Item(5, 18).Value = "Hortaliza"
